$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: Accelerometer, $9.95, hyperlink to SparkFun product page
$ws.Range("A5").Value = "Accelerometer"

$ws.Range("B5").Value = 9.95
$ws.Range("B5").NumberFormat = $ws.Range("B3").NumberFormat

$ws.Range("C5").Value = "https://www.sparkfun.com/products/12786"
$ws.Hyperlinks.Add($ws.Range("C5"), "https://www.sparkfun.com/products/12786") | Out-Null
$ws.Range("C5").Style = "Hyperlink"

# Widen column A to fit the new, longer labels
$ws.Columns.Item(1).ColumnWidth = 19

# Match the saved selection/cursor position
$ws.Range("G3").Select()
